$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (rows 2-7), replacing old rows 2-11
$data = @(
    @(0, 2, "SMART SENSING MIDDLEWARE", 78.20999999999999, 1, 78.20999999999999, 2),
    @(1, 2, "SHAMIYANA APP", 64.5, 0.85, 54.82, 2),
    @(3, 3, "Post-processing of Large Language Models", 111.25, 1, 100, 4),
    @(4, 3, "Multi Model Data Analysis for Annotation of Human Activities", 104.4, 1, 100, 4),
    @(5, 3, "Cloudphysician's Vital Extraction Challenge", 78.20999999999999, 0.8, 62.57, 4),
    @(6, 3, "Website for the Literature Society of the college", 74.88, 0.85, 63.65, 4)
)

# Clear old rows 2-11 first (old sheet had rows up to 11).
# Use ClearContents on rows that remain (2-7) to preserve column A's style,
# and Clear (contents+format) on rows that are being fully removed (8-11).
$ws.Range("A2:G7").ClearContents()
$ws.Range("A8:G11").Clear()

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r++
}
